$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# The c0017c90 file went through a newer Handback -> Xliff generation pass;
# it is still "in sync with en-US" but the generation timestamp moved forward.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2017-02-09 06:18:19"
# Column got a bit wider to fit the new (longer) status strings used on the
# language sheets.
$wsOverview.Range("E1:F1").ColumnWidth = 32.59

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# 1088303c file fell out of sync with en-US.
$wsZhCn.Range("C2").Value = "Handed back: not in sync with en-US"
# c0017c90 file went through a new handoff/handback round.
$wsZhCn.Range("H3").Value = "2017-02-09 06:18:02"
$wsZhCn.Range("L3").Value = "2017-02-09 06:18:57"
# Status column widened to fit the longer "not in sync" text.
$wsZhCn.Range("C1").ColumnWidth = 32.59

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# 1088303c file fell out of sync with en-US.
$wsDeDe.Range("C2").Value = "Handed back: not in sync with en-US"
# c0017c90 file went through a new handoff/handback round.
$wsDeDe.Range("H3").Value = "2017-02-09 06:18:19"
$wsDeDe.Range("L3").Value = "2017-02-09 06:19:22"
# Status column widened to fit the longer "not in sync" text.
$wsDeDe.Range("C1").ColumnWidth = 32.59
